# Remove the "SAE outcome" (AE_AESER_AEOUT) row from the common_forms sheet.
# This field is not populated in clinsight metadata, so the whole metadata
# row describing it is deleted (rather than e.g. blanking individual cells).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("common_forms")

# Row 17 holds: A17 = "AE_AESER_AEOUT", C17 = "SAE outcome" (D17 = "other",
# E17 = "Adverse events"). Deleting the entire row shifts the rows below it
# up by one, shrinks the table/sheet dimensions accordingly, and drops the
# now-unused shared strings automatically.
$ws.Rows.Item(17).Delete()
